$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD2 = $ws.Range("D2")
$sD2 = $cD2.Style
$cD2.Value = "'63.879.11"
$cD2.Style = $sD2
$ws.Range("E2").Value = "  +2.82%  "

$cD3 = $ws.Range("D3")
$sD3 = $cD3.Style
$cD3.Value = "'2.542.49"
$cD3.Style = $sD3
$ws.Range("E3").Value = "  +5.18%  "

$ws.Range("E4").Value = "  -0.06%  "

$cD5 = $ws.Range("D5")
$sD5 = $cD5.Style
$cD5.Value = "'574.01"
$cD5.Style = $sD5
$ws.Range("E5").Value = "  +2.20%  "

$cD6 = $ws.Range("D6")
$sD6 = $cD6.Style
$cD6.Value = "'148.50"
$cD6.Style = $sD6
$ws.Range("E6").Value = "  +6.86%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  +0.43%  "

$cD9 = $ws.Range("D9")
$sD9 = $cD9.Style
$cD9.Value = "'2.541.59"
$cD9.Style = $sD9
$ws.Range("E9").Value = "  +5.29%  "

$ws.Range("E10").Value = "  +2.43%  "

$cD11 = $ws.Range("D11")
$sD11 = $cD11.Style
$cD11.Value = "'5.79"
$cD11.Style = $sD11
$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("E12").Value = "  +1.63%  "

$cD13 = $ws.Range("D13")
$sD13 = $cD13.Style
$cD13.Value = "'0.361"
$cD13.Style = $sD13
$ws.Range("E13").Value = "  +3.23%  "

$cD14 = $ws.Range("D14")
$sD14 = $cD14.Style
$cD14.Value = "'28.05"
$cD14.Style = $sD14
$ws.Range("E14").Value = "  +8.69%  "

$cD15 = $ws.Range("D15")
$sD15 = $cD15.Style
$cD15.Value = "'2.995.74"
$cD15.Style = $sD15
$ws.Range("E15").Value = "  +5.19%  "

$cD16 = $ws.Range("D16")
$sD16 = $cD16.Style
$cD16.Value = "'63.653.54"
$cD16.Style = $sD16
$ws.Range("E16").Value = "  +2.58%  "

$ws.Range("E17").Value = "  +3.75%  "

$cD18 = $ws.Range("D18")
$sD18 = $cD18.Style
$cD18.Value = "'2.540.29"
$cD18.Style = $sD18
$ws.Range("E18").Value = "  +5.30%  "

$cD19 = $ws.Range("D19")
$sD19 = $cD19.Style
$cD19.Value = "'11.61"
$cD19.Style = $sD19
$ws.Range("E19").Value = "  +4.80%  "

$cD20 = $ws.Range("D20")
$sD20 = $cD20.Style
$cD20.Value = "'344.52"
$cD20.Style = $sD20
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("E21").Value = "  +3.49%  "

$cD22 = $ws.Range("D22")
$sD22 = $cD22.Style
$cD22.Value = "'6.89"
$cD22.Style = $sD22
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("E23").Value = "  +0.53%  "

$cD24 = $ws.Range("D24")
$sD24 = $cD24.Style
$cD24.Value = "'66.17"
$cD24.Style = $sD24
$ws.Range("E24").Value = "  +1.82%  "

$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E26").Value = "  +5.10%  "

$ws.Range("E27").Value = "  +0.21%  "

$cD28 = $ws.Range("D28")
$sD28 = $cD28.Style
$cD28.Value = "'8.33"
$cD28.Style = $sD28
$ws.Range("E28").Value = "  -0.60%  "

$cD29 = $ws.Range("D29")
$sD29 = $cD29.Style
$cD29.Value = "'1.42"
$cD29.Style = $sD29
$ws.Range("E29").Value = "  +3.60%  "

$cD30 = $ws.Range("D30")
$sD30 = $cD30.Style
$cD30.Value = "'0.0₃0828"
$cD30.Style = $sD30
$ws.Range("E30").Value = "  +6.68%  "

$ws.Range("E31").Value = "  +4.36%  "

$cD32 = $ws.Range("D32")
$sD32 = $cD32.Style
$cD32.Value = "'6.83"
$cD32.Style = $sD32
$ws.Range("E32").Value = "  +6.97%  "

$cD33 = $ws.Range("D33")
$sD33 = $cD33.Style
$cD33.Value = "'176.81"
$cD33.Style = $sD33
$ws.Range("E33").Value = "  +3.62%  "

$cD34 = $ws.Range("D34")
$sD34 = $cD34.Style
$cD34.Value = "'1.58"
$cD34.Style = $sD34
$ws.Range("E34").Value = "  +12.76%  "

$cD35 = $ws.Range("D35")
$sD35 = $cD35.Style
$cD35.Value = "'424.93"
$cD35.Style = $sD35
$ws.Range("E35").Value = "  +16.49%  "

$ws.Range("E36").Value = "  +3.22%  "

$cD37 = $ws.Range("D37")
$sD37 = $cD37.Style
$cD37.Value = "'19.16"
$cD37.Style = $sD37
$ws.Range("E37").Value = "  +3.04%  "

$cD38 = $ws.Range("D38")
$sD38 = $cD38.Style
$cD38.Value = "'4.46"
$cD38.Style = $sD38
$ws.Range("E38").Value = "  -1.99%  "

$ws.Range("E39").Value = "  +0.00%  "

$cD40 = $ws.Range("D40")
$sD40 = $cD40.Style
$cD40.Value = "'1.76"
$cD40.Style = $sD40
$ws.Range("E40").Value = "  +5.26%  "

$ws.Range("E41").Value = "  -0.02%  "

$cD42 = $ws.Range("D42")
$sD42 = $cD42.Style
$cD42.Value = "'40.59"
$cD42.Style = $sD42
$ws.Range("E42").Value = "  +3.88%  "

$cD43 = $ws.Range("D43")
$sD43 = $cD43.Style
$cD43.Value = "'152.74"
$cD43.Style = $sD43
$ws.Range("E43").Value = "  +5.89%  "

$ws.Range("E44").Value = "  +3.59%  "

$ws.Range("E45").Value = "  +2.42%  "

$cD46 = $ws.Range("D46")
$sD46 = $cD46.Style
$cD46.Value = "'0.613"
$cD46.Style = $sD46
$ws.Range("E46").Value = "  +4.69%  "

$cD47 = $ws.Range("D47")
$sD47 = $cD47.Style
$cD47.Value = "'0.0534"
$cD47.Style = $sD47
$ws.Range("E47").Value = "  +2.80%  "

$cD48 = $ws.Range("D48")
$sD48 = $cD48.Style
$cD48.Value = "'0.0969"
$cD48.Style = $sD48

$cD49 = $ws.Range("D49")
$sD49 = $cD49.Style
$cD49.Value = "'18.96"
$cD49.Style = $sD49
$ws.Range("E49").Value = "  +5.42%  "

$ws.Range("E50").Value = "  +5.08%  "

$cD51 = $ws.Range("D51")
$sD51 = $cD51.Style
$cD51.Value = "'1.84"
$cD51.Style = $sD51
$ws.Range("E51").Value = "  +8.78%  "
